# Refresh the crypto price/volume table (columns D and E, rows 2-51)
# with updated values, matching the upstream GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.878.32"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").Value = "1.729.50"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9974"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.85"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9978"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4910"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2603"
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06221"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "1.734.93"
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.05"
$ws.Range("E11").Value = "  +3.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06911"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6104"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.501"
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.39"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9983"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "26.637.91"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9974"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007195"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").Value = "1.953.60"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.434"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.569"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.82"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.35"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.787"
$ws.Range("E27").Value = "  +5.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.382"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.38"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.956"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08003"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.688"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04534"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9971"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.609"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6246"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9341"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.057"
$ws.Range("E39").Value = "  +5.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.444"
$ws.Range("E40").Value = "  +2.19%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01503"
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.664"
$ws.Range("E43").Value = "  +4.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.78"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3868"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.947"
$ws.Range("E46").Value = "  +3.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1163"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05390"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.921"
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.27"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.242"
$ws.Range("E51").Value = "  +0.07%  "
